# The test E-Mail addresses embedded in the "UsuariosRegistro" sheet (and
# reused by the "LoginData" sheet, since both point at the same shared
# strings) carry a generated timestamp as part of the address, e.g.
#   juan.perez+20251109_011412@test.com
# Refresh the stamp to 20251109_012452 everywhere it appears in the
# workbook so every worksheet referencing that text stays in sync.

$wb = $excel.ActiveWorkbook

$oldStamp = "20251109_011412"
$newStamp = "20251109_012452"

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldStamp, $newStamp)
}
